$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1840
$ws.Range("F8").Value = 13222
$ws.Range("F9").Value = 13084
$ws.Range("F12").Value = 17
$ws.Range("F19").Value = 47
$ws.Range("F21").Value = 209
$ws.Range("F23").Value = 742

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 15

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 14

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1840
$ws.Range("F10").Value = 13222
$ws.Range("F11").Value = 13084
$ws.Range("F14").Value = 17
$ws.Range("F23").Value = 47
$ws.Range("F27").Value = 14
$ws.Range("F28").Value = 209
$ws.Range("F30").Value = 742
$ws.Range("F33").Value = 15
